$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 352, shifting existing rows 352..363 down to 353..364.
$ws.Rows.Item(352).Insert()

# Populate the newly inserted row 352 with the new weekly record.
# (Same shape as the surrounding "Feria Lagunitas de Puerto Montt" / Ciboulette rows,
# copying row 353's non-changing columns.)
$ws.Range("A352").Value = 4
$ws.Range("B352").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C352").Value = "Los Lagos"
$ws.Range("D352").Value = 45075
$ws.Range("E352").Value = 10
$ws.Range("F352").Value = 100112039
$ws.Range("G352").Value = "Ciboulette"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 80
$ws.Range("K352").Value = 3500
$ws.Range("L352").Value = 3500
$ws.Range("M352").Value = 3500
$ws.Range("N352").Value = "$/docena de atados"
$ws.Range("O352").Value = "Región Metropolitana"
$ws.Range("P352").Value = 1167
$ws.Range("Q352").Value = 3
$ws.Range("R352").Value = "Hortaliza"
